# Apply BOM update: add new "Akku" (battery) line item as row 26,
# shifting the trailing SUM row from 28 down to 29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26 (this pushes the blank rows 26/27 and the
# SUM row 28 down by one, so SUM ends up on row 29 - matching the target).
$ws.Rows.Item(26).Insert()

# Fill in the new row 26 with the battery BOM line.
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 3).Value = "Akku"
$ws.Cells.Item(26, 4).Value = 2
$ws.Cells.Item(26, 6).Value = 12.95
$ws.Cells.Item(26, 7).Formula = "=F26*D26"

# Add the hyperlink in column H for row 26 (visible text = the URL itself,
# matching the convention used by the other BOM link cells in this sheet).
$akkuUrl = "https://www.conrad.ch/de/p/samsung-inr18650-29e-spezial-akku-18650-flat-top-hochtemperaturfaehig-li-ion-3-6-v-2900-mah-2239492.html?searchType=SearchRedirect"
$ws.Hyperlinks.Add($ws.Cells.Item(26, 8), $akkuUrl)
$ws.Cells.Item(26, 8).Style = "Hyperlink"

# Update the trailing SUM formula (now on row 29) to include the new row.
$ws.Cells.Item(29, 7).Formula = "=SUM(G2:G26)"

# Update the selection to match the target state.
$ws.Range("H27").Select()

Write-Host "Done"
